# Applies the OOXML changes described by the diff:
#  - Sez Urto: scroll/selection moves to G127 (topLeftCell A57), cell H126 (empty, styled) removed
#  - Decad K: scroll resets to top (no explicit topLeftCell)
#  - Fabrizio syst: becomes the active/selected tab
#  - Giorgio syst: loses tabSelected, zoom resets to 100, selection moves to I20,
#                  new column B ("Anti-p entranti") with data for rows 2-15,
#                  and corrected counts in rows 19-27
#  - workbook: activeTab points at Fabrizio syst (index 2)

$wb = $excel.ActiveWorkbook

# ---- Sheet "Sez Urto" ----
$wsSezUrto = $wb.Worksheets.Item("Sez Urto")
$wsSezUrto.Activate()
$excel.ActiveWindow.ScrollRow = 57
$excel.ActiveWindow.ScrollColumn = 1
$wsSezUrto.Range("G127").Select() | Out-Null
$wsSezUrto.Range("H126").Clear()

# ---- Sheet "Decad K" ----
$wsDecadK = $wb.Worksheets.Item("Decad K")
$wsDecadK.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsDecadK.Range("D19").Select() | Out-Null

# ---- Sheet "Giorgio syst" ----
$wsGiorgio = $wb.Worksheets.Item("Giorgio syst")
$wsGiorgio.Activate()
$excel.ActiveWindow.Zoom = 100

$wsGiorgio.Columns.Item(2).ColumnWidth = 21.14

$wsGiorgio.Range("B1").Value = "Anti-p entranti"

$giorgioData = @(
  @(14, 3, 2),
  @(6, 1, 2),
  @(13, 2, 1),
  @(6, 0, 2),
  @(11, 2, 3),
  @(6, 2, 0),
  @(8, 1, 2),
  @(5, 1, 2),
  @(10, 2, 3),
  @(9, 0, 2),
  @(11, 0, 2),
  @(8, 1, 1),
  @(12, 1, 4),
  @(11, 2, 2)
)
for ($i = 0; $i -lt $giorgioData.Length; $i++) {
  $r = 2 + $i
  $wsGiorgio.Cells.Item($r, 2).Value = $giorgioData[$i][0]
  $wsGiorgio.Cells.Item($r, 3).Value = $giorgioData[$i][1]
  $wsGiorgio.Cells.Item($r, 4).Value = $giorgioData[$i][2]
}

$giorgioCounts = @{
  19 = @(14, 0)
  20 = @(18, 1)
  21 = @(14, 1)
  22 = @(14, 1)
  23 = @(12, 0)
  24 = @(10, 1)
  25 = @(12, 1)
  26 = @(10, 2)
  27 = @(9, 0)
}
foreach ($r in $giorgioCounts.Keys) {
  $vals = $giorgioCounts[$r]
  $wsGiorgio.Cells.Item($r, 2).Value = $vals[0]
  $wsGiorgio.Cells.Item($r, 3).Value = $vals[1]
}

$wsGiorgio.Range("I20").Select() | Out-Null

# ---- Sheet "Fabrizio syst" becomes the active tab ----
$wsFabrizio = $wb.Worksheets.Item("Fabrizio syst")
$wsFabrizio.Activate()
